$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row in the sheet at row 151 (pushes existing rows 151+ down by one)
$ws.Rows.Item(151).Insert()

# Copy formatting from the row above (Arial font style, wrapped text in column E)
# onto the newly-inserted blank row so it matches its table neighbours.
$ws.Range("A150:E150").Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)

# Expand the table (ListObject) to include the newly inserted row
$lo.Resize($ws.Range("A1:E291"))

# Populate the new row with the insurance derived variable
$ws.Range("A151").Value = "D20"
$ws.Range("B151").Value = "insurance"
$ws.Range("C151").Value = "Demographics"
$ws.Range("D151").Value = "Insurance type"
$ws.Range("E151").Value = "Medicaid alone; Medicare alone; Medicare/Medicaid +/- other; Other government +/- other; Private +/- other; Uninsured; Unknown"

# Row height for the wrapped, multi-line Values text
$ws.Rows.Item(151).RowHeight = 46

# Reflect the cell the author ended up viewing/selecting after the edit
$ws.Range("E151").Select()
